# Auto-generated: update cryptos price (D) and volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Values are forced to text via a leading apostrophe (Excel "quote prefix")
# so strings that look numeric ("1.002", "218.00", ...) are not silently
# converted to numbers/rounded, matching the original inline-string cells.
# The Style reset afterwards clears the quote-prefix formatting flag so the
# cell keeps its original (unstyled) appearance.
$priceUpdates = @{
    2 = '26.043.84'
    3 = '1.646.86'
    4 = '1.002'
    5 = '218.00'
    6 = '0.5183'
    7 = '1.003'
    8 = '0.2620'
    9 = '0.06292'
    10 = '20.30'
    11 = '0.07664'
    12 = '4.573'
    13 = '1.651.02'
    14 = '1.872.48'
    15 = '0.5563'
    16 = '0.0₅8105'
    17 = '65.00'
    18 = '26.007.12'
    20 = '4.594'
    21 = '10.42'
    22 = '191.82'
    23 = '5.893'
    24 = '1.003'
    25 = '144.16'
    27 = '7.165'
    28 = '15.80'
    29 = '1.514'
    30 = '0.05339'
    32 = '3.441'
    33 = '3.330'
    34 = '1.545'
    35 = '2.416'
    36 = '2.775'
    37 = '0.9396'
    38 = '0.5580'
    39 = '0.01571'
    40 = '5.768'
    41 = '1.004'
    42 = '1.029.90'
    43 = '0.8237'
    44 = '100.69'
    45 = '1.782.61'
    47 = '57.05'
    48 = '0.9991'
    49 = '0.4310'
    50 = '7.881'
}
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.Value = "'" + $priceUpdates[$row]
    $cell.Style = "Normal"
}

# --- Volume(1h) (column E) updates ---
# These already carry padding spaces and a trailing "%" so Excel stores them
# as plain text without any extra handling.
$volumeUpdates = @{
    2 = '  -0.04%  '
    3 = '  +0.21%  '
    4 = '  -0.17%  '
    5 = '  +0.19%  '
    6 = '  -0.04%  '
    7 = '  -0.14%  '
    8 = '  +0.26%  '
    9 = '  +0.27%  '
    10 = '  -0.51%  '
    11 = '  -1.16%  '
    12 = '  +2.37%  '
    13 = '  +2.64%  '
    14 = '  +0.15%  '
    15 = '  -0.20%  '
    16 = '  +1.55%  '
    17 = '  +0.47%  '
    18 = '  -0.19%  '
    19 = '  -0.11%  '
    20 = '  -0.61%  '
    21 = '  +3.47%  '
    22 = '  -0.32%  '
    23 = '  -0.87%  '
    24 = '  -0.23%  '
    26 = '  -1.63%  '
    27 = '  +0.10%  '
    28 = '  -0.44%  '
    29 = '  +2.28%  '
    30 = '  -5.14%  '
    31 = '  +0.19%  '
    32 = '  -0.40%  '
    33 = '  -0.60%  '
    34 = '  -3.12%  '
    35 = '  +0.23%  '
    36 = '  -0.48%  '
    37 = '  +0.38%  '
    38 = '  -1.36%  '
    39 = '  -0.18%  '
    40 = '  -2.99%  '
    42 = '  -1.99%  '
    43 = '  -1.93%  '
    44 = '  -1.52%  '
    45 = '  +0.11%  '
    46 = '  +10.04%  '
    47 = '  +0.35%  '
    48 = '  -0.94%  '
    49 = '  -0.35%  '
    50 = '  -0.36%  '
    51 = '  -3.76%  '
}
foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

